$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.3719002097481816
$ws.Range("F3").Value = 1.3760887635139747
$ws.Range("F4").Value = 1.3642754662155234
$ws.Range("F5").Value = 1.3683573211496605
$ws.Range("F6").Value = 1.3673349283339424
$ws.Range("F7").Value = 1.3648672142025819
$ws.Range("F8").Value = 1.0685113476047834
$ws.Range("F9").Value = 1.0835459285056379
$ws.Range("F10").Value = 1.0552534968012119
$ws.Range("F11").Value = 0.14264426372740593
$ws.Range("F12").Value = 0.16753990864725479
$ws.Range("F13").Value = 0.1618914352032011

$ws.Range("H18").Select()
